$wb = $excel.ActiveWorkbook

# Gungnir_Profits data refresh (scheduled runner) -- updates computed
# market-price / profit columns (H:N) across all job sheets.

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 9260853  # H2
$ws.Cells.Item(2, 9).Value = 20833384  # I2
$ws.Cells.Item(2, 10).Value = 2827.8  # J2
$ws.Cells.Item(2, 11).Value = 20833384  # K2
$ws.Cells.Item(2, 12).Value = 2827.8  # L2
$ws.Cells.Item(2, 13).Value = -20833271  # M2
$ws.Cells.Item(2, 14).Value = -3053.8  # N2
$ws.Cells.Item(100, 8).Value = 5309.5557  # H100
$ws.Cells.Item(100, 9).Value = 6547.6665  # I100
$ws.Cells.Item(100, 11).Value = 6547.6665  # K100
$ws.Cells.Item(100, 13).Value = -6006.6665  # M100
$ws.Cells.Item(132, 8).Value = 3189.75  # H132
$ws.Cells.Item(132, 9).Value = 2924.233  # I132
$ws.Cells.Item(132, 10).Value = 3907.6296  # J132
$ws.Cells.Item(132, 11).Value = 8772.699000000001  # K132
$ws.Cells.Item(132, 12).Value = 11722.8888  # L132
$ws.Cells.Item(132, 13).Value = -6242.699000000001  # M132
$ws.Cells.Item(132, 14).Value = -16782.8888  # N132
$ws.Cells.Item(137, 8).Value = 898.7458  # H137
$ws.Cells.Item(137, 9).Value = 788.5599999999999  # I137
$ws.Cells.Item(137, 11).Value = 2365.68  # K137
$ws.Cells.Item(137, 13).Value = 184.3200000000002  # M137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1070  # H2
$ws.Cells.Item(2, 9).Value = 1023.3333  # I2
$ws.Cells.Item(2, 10).Value = 1163.3334  # J2
$ws.Cells.Item(2, 11).Value = 1023.3333  # K2
$ws.Cells.Item(2, 12).Value = 1163.3334  # L2
$ws.Cells.Item(2, 13).Value = -910.3333  # M2
$ws.Cells.Item(2, 14).Value = -1389.3334  # N2
$ws.Cells.Item(50, 8).Value = 1785.8334  # H50
$ws.Cells.Item(50, 9).Value = 349.33334  # I50
$ws.Cells.Item(50, 10).Value = 2264.6667  # J50
$ws.Cells.Item(50, 11).Value = 349.33334  # K50
$ws.Cells.Item(50, 12).Value = 2264.6667  # L50
$ws.Cells.Item(50, 13).Value = 364.66666  # M50
$ws.Cells.Item(50, 14).Value = -3692.6667  # N50
$ws.Cells.Item(74, 8).Value = 1184.909  # H74
$ws.Cells.Item(74, 9).Value = 1167.8064  # I74
$ws.Cells.Item(74, 10).Value = 1450  # J74
$ws.Cells.Item(74, 11).Value = 1167.8064  # K74
$ws.Cells.Item(74, 12).Value = 1450  # L74
$ws.Cells.Item(74, 13).Value = -293.8063999999999  # M74
$ws.Cells.Item(74, 14).Value = -3198  # N74
$ws.Cells.Item(77, 8).Value = 1184.909  # H77
$ws.Cells.Item(77, 9).Value = 1167.8064  # I77
$ws.Cells.Item(77, 10).Value = 1450  # J77
$ws.Cells.Item(77, 11).Value = 5839.031999999999  # K77
$ws.Cells.Item(77, 12).Value = 7250  # L77
$ws.Cells.Item(77, 13).Value = -1471.031999999999  # M77
$ws.Cells.Item(77, 14).Value = -15986  # N77
$ws.Cells.Item(116, 8).Value = 1070  # H116
$ws.Cells.Item(116, 9).Value = 1023.3333  # I116
$ws.Cells.Item(116, 10).Value = 1163.3334  # J116
$ws.Cells.Item(116, 11).Value = 1023.3333  # K116
$ws.Cells.Item(116, 12).Value = 1163.3334  # L116
$ws.Cells.Item(116, 13).Value = 1270.6667  # M116
$ws.Cells.Item(116, 14).Value = -5751.3334  # N116
$ws.Cells.Item(132, 8).Value = 1252796.5  # H132
$ws.Cells.Item(132, 9).Value = 900.4146  # I132
$ws.Cells.Item(132, 10).Value = 9807419  # J132
$ws.Cells.Item(132, 11).Value = 2701.2438  # K132
$ws.Cells.Item(132, 12).Value = 29422257  # L132
$ws.Cells.Item(132, 13).Value = -171.2437999999997  # M132
$ws.Cells.Item(132, 14).Value = -29427317  # N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1070  # H3
$ws.Cells.Item(3, 9).Value = 1023.3333  # I3
$ws.Cells.Item(3, 10).Value = 1163.3334  # J3
$ws.Cells.Item(3, 11).Value = 1023.3333  # K3
$ws.Cells.Item(3, 12).Value = 1163.3334  # L3
$ws.Cells.Item(3, 13).Value = -909.3333  # M3
$ws.Cells.Item(3, 14).Value = -1391.3334  # N3
$ws.Cells.Item(20, 8).Value = 4074.3333  # H20
$ws.Cells.Item(20, 10).Value = 4029  # J20
$ws.Cells.Item(20, 12).Value = 4029  # L20
$ws.Cells.Item(20, 14).Value = -4523  # N20
$ws.Cells.Item(105, 8).Value = 45456120  # H105
$ws.Cells.Item(105, 9).Value = 1521.0555  # I105
$ws.Cells.Item(105, 10).Value = 250001820  # J105
$ws.Cells.Item(105, 11).Value = 1521.0555  # K105
$ws.Cells.Item(105, 12).Value = 250001820  # L105
$ws.Cells.Item(105, 13).Value = 225.9445000000001  # M105
$ws.Cells.Item(105, 14).Value = -250005314  # N105
$ws.Cells.Item(134, 8).Value = 3631.976  # H134
$ws.Cells.Item(134, 9).Value = 898.35297  # I134
$ws.Cells.Item(134, 11).Value = 2695.05891  # K134
$ws.Cells.Item(134, 13).Value = -160.0589100000002  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1384.1698  # H31
$ws.Cells.Item(31, 9).Value = 1084.2222  # I31
$ws.Cells.Item(31, 10).Value = 1538.4286  # J31
$ws.Cells.Item(31, 11).Value = 1084.2222  # K31
$ws.Cells.Item(31, 12).Value = 1538.4286  # L31
$ws.Cells.Item(31, 13).Value = -789.2221999999999  # M31
$ws.Cells.Item(31, 14).Value = -2128.4286  # N31
$ws.Cells.Item(34, 8).Value = 1384.1698  # H34
$ws.Cells.Item(34, 9).Value = 1084.2222  # I34
$ws.Cells.Item(34, 10).Value = 1538.4286  # J34
$ws.Cells.Item(34, 11).Value = 1084.2222  # K34
$ws.Cells.Item(34, 12).Value = 1538.4286  # L34
$ws.Cells.Item(34, 13).Value = -882.2221999999999  # M34
$ws.Cells.Item(34, 14).Value = -1942.4286  # N34
$ws.Cells.Item(35, 8).Value = 2365.4666  # H35
$ws.Cells.Item(35, 9).Value = 2280.923  # I35
$ws.Cells.Item(35, 10).Value = 2915  # J35
$ws.Cells.Item(35, 11).Value = 2280.923  # K35
$ws.Cells.Item(35, 12).Value = 2915  # L35
$ws.Cells.Item(35, 13).Value = -1986.923  # M35
$ws.Cells.Item(35, 14).Value = -3503  # N35
$ws.Cells.Item(58, 8).Value = 18519500  # H58
$ws.Cells.Item(58, 9).Value = 27778814  # I58
$ws.Cells.Item(58, 10).Value = 874.55554  # J58
$ws.Cells.Item(58, 11).Value = 27778814  # K58
$ws.Cells.Item(58, 12).Value = 874.55554  # L58
$ws.Cells.Item(58, 13).Value = -27778611  # M58
$ws.Cells.Item(58, 14).Value = -1280.55554  # N58
$ws.Cells.Item(132, 8).Value = 5377413  # H132
$ws.Cells.Item(132, 9).Value = 754.7273  # I132
$ws.Cells.Item(132, 10).Value = 18520356  # J132
$ws.Cells.Item(132, 11).Value = 2264.1819  # K132
$ws.Cells.Item(132, 12).Value = 55561068  # L132
$ws.Cells.Item(132, 13).Value = 265.8181  # M132
$ws.Cells.Item(132, 14).Value = -55566128  # N132
$ws.Cells.Item(136, 8).Value = 18519500  # H136
$ws.Cells.Item(136, 9).Value = 27778814  # I136
$ws.Cells.Item(136, 10).Value = 874.55554  # J136
$ws.Cells.Item(136, 11).Value = 83336442  # K136
$ws.Cells.Item(136, 12).Value = 2623.66662  # L136
$ws.Cells.Item(136, 13).Value = -83333892  # M136
$ws.Cells.Item(136, 14).Value = -7723.66662  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 25645428  # H5
$ws.Cells.Item(5, 9).Value = 39216056  # I5
$ws.Cells.Item(5, 10).Value = 12021.444  # J5
$ws.Cells.Item(5, 11).Value = 117648168  # K5
$ws.Cells.Item(5, 12).Value = 36064.33199999999  # L5
$ws.Cells.Item(5, 13).Value = -117648056  # M5
$ws.Cells.Item(5, 14).Value = -36288.33199999999  # N5
$ws.Cells.Item(12, 8).Value = 22465.756  # H12
$ws.Cells.Item(12, 10).Value = 29751.027  # J12
$ws.Cells.Item(12, 12).Value = 89253.08099999999  # L12
$ws.Cells.Item(12, 14).Value = -89599.08099999999  # N12
$ws.Cells.Item(132, 8).Value = 29418730  # H132
$ws.Cells.Item(132, 9).Value = 842.6  # I132
$ws.Cells.Item(132, 10).Value = 71444290  # J132
$ws.Cells.Item(132, 11).Value = 7583.400000000001  # K132
$ws.Cells.Item(132, 12).Value = 642998610  # L132
$ws.Cells.Item(132, 13).Value = -5053.400000000001  # M132
$ws.Cells.Item(132, 14).Value = -643003670  # N132
$ws.Cells.Item(133, 8).Value = 83338300  # H133
$ws.Cells.Item(133, 9).Value = 111116130  # I133
$ws.Cells.Item(133, 11).Value = 333348390  # K133
$ws.Cells.Item(133, 13).Value = -333343330  # M133
$ws.Cells.Item(135, 8).Value = 25645428  # H135
$ws.Cells.Item(135, 9).Value = 39216056  # I135
$ws.Cells.Item(135, 10).Value = 12021.444  # J135
$ws.Cells.Item(135, 11).Value = 352944504  # K135
$ws.Cells.Item(135, 12).Value = 108192.996  # L135
$ws.Cells.Item(135, 13).Value = -352941969  # M135
$ws.Cells.Item(135, 14).Value = -113262.996  # N135
$ws.Cells.Item(139, 8).Value = 222775.8  # H139
$ws.Cells.Item(139, 9).Value = 1065.6428  # I139
$ws.Cells.Item(139, 10).Value = 533170  # J139
$ws.Cells.Item(139, 11).Value = 3196.9284  # K139
$ws.Cells.Item(139, 12).Value = 1599510  # L139
$ws.Cells.Item(139, 13).Value = 1943.0716  # M139
$ws.Cells.Item(139, 14).Value = -1609790  # N139

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 12071.523  # H132
$ws.Cells.Item(132, 9).Value = 9283  # I132
$ws.Cells.Item(132, 10).Value = 15789.556  # J132
$ws.Cells.Item(132, 11).Value = 27849  # K132
$ws.Cells.Item(132, 12).Value = 47368.66800000001  # L132
$ws.Cells.Item(132, 13).Value = -25319  # M132
$ws.Cells.Item(132, 14).Value = -52428.66800000001  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 822.34283  # H22
$ws.Cells.Item(22, 9).Value = 383.46155  # I22
$ws.Cells.Item(22, 10).Value = 1081.6818  # J22
$ws.Cells.Item(22, 11).Value = 383.46155  # K22
$ws.Cells.Item(22, 12).Value = 1081.6818  # L22
$ws.Cells.Item(22, 13).Value = -88.46154999999999  # M22
$ws.Cells.Item(22, 14).Value = -1671.6818  # N22
$ws.Cells.Item(27, 8).Value = 822.34283  # H27
$ws.Cells.Item(27, 9).Value = 383.46155  # I27
$ws.Cells.Item(27, 10).Value = 1081.6818  # J27
$ws.Cells.Item(27, 11).Value = 383.46155  # K27
$ws.Cells.Item(27, 12).Value = 1081.6818  # L27
$ws.Cells.Item(27, 13).Value = -276.46155  # M27
$ws.Cells.Item(27, 14).Value = -1295.6818  # N27
$ws.Cells.Item(132, 8).Value = 16398665  # H132
$ws.Cells.Item(132, 9).Value = 32259668  # I132
$ws.Cells.Item(132, 10).Value = 8961.032999999999  # J132
$ws.Cells.Item(132, 11).Value = 96779004  # K132
$ws.Cells.Item(132, 12).Value = 26883.099  # L132
$ws.Cells.Item(132, 13).Value = -96776474  # M132
$ws.Cells.Item(132, 14).Value = -31943.099  # N132
$ws.Cells.Item(136, 8).Value = 32655470  # H136
$ws.Cells.Item(136, 9).Value = 5104709  # I136
$ws.Cells.Item(136, 10).Value = 142858510  # J136
$ws.Cells.Item(136, 11).Value = 15314127  # K136
$ws.Cells.Item(136, 12).Value = 428575530  # L136
$ws.Cells.Item(136, 13).Value = -15311577  # M136
$ws.Cells.Item(136, 14).Value = -428580630  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 16280.747  # H132
$ws.Cells.Item(132, 9).Value = 21067.04  # I132
$ws.Cells.Item(132, 11).Value = 63201.12  # K132
$ws.Cells.Item(132, 13).Value = -60671.12  # M132
